$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.89402836561203
$ws.Range("B1").Value = 1.781641840934753
$ws.Range("C1").Value = 4.17420482635498
$ws.Range("D1").Value = 3.497797966003418
$ws.Range("E1").Value = 1.507249593734741
